$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new "EMAIL" column before the existing GENDER column (C).
#    Everything from the old column C onward shifts right by one (C->D,
#    D->E, ... K->L), which also moves the B/C style-only filler rows
#    (8-18) and column width definitions along with it, and preserves the
#    date-formatted style already sitting in (the now-shifted) column F.
# ---------------------------------------------------------------------------
$ws.Columns("C").Insert()
$ws.Range("C1").Value = "EMAIL"
$ws.Columns("C").ColumnWidth = 14

# ---------------------------------------------------------------------------
# 2. Update the (already existing) first data row with the new mobile
#    number / DOB, reusing its own short-date number format for the newly
#    added rows below it.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = 9876543210
$ws.Range("F2").Copy()
$ws.Range("F3:F6").PasteSpecial(-4122)
$ws.Range("F2").Value = 32880

# ---------------------------------------------------------------------------
# 3. Add the four new people (rows 3-6).
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Karthick"
$ws.Range("B3").Value = "Subramaniam"
$ws.Range("D3").Value = "Male"
$ws.Range("E3").Value = 9876543210
$ws.Range("F3").Value = 33277
$ws.Range("G3").Value = "MATHS, PHYSICS"
$ws.Range("H3").Value = "SPORTS"
$ws.Range("I3").Value = "GUINDY, CHENNAI"
$ws.Range("J3").Value = "Haryana"
$ws.Range("K3").Value = "Panipat"

$ws.Range("A4").Value = "Saranya"
$ws.Range("B4").Value = "Santhosh"
$ws.Range("D4").Value = "Female"
$ws.Range("E4").Value = 9876543210
$ws.Range("F4").Value = 33672
$ws.Range("G4").Value = "TAMIL, MATHS"
$ws.Range("H4").Value = "READING, MUSIC"
$ws.Range("I4").Value = "BANGLORE"
$ws.Range("J4").Value = "Haryana"
$ws.Range("K4").Value = "Panipat"

$ws.Range("A5").Value = "Parthiban"
$ws.Range("B5").Value = "Sankaran"
$ws.Range("D5").Value = "Male"
$ws.Range("E5").Value = 9876543210
$ws.Range("F5").Value = 34069
$ws.Range("G5").Value = "ENGLISH, PHYSICS"
$ws.Range("H5").Value = "MUSIC"
$ws.Range("I5").Value = "SALEM"
$ws.Range("J5").Value = "Haryana"
$ws.Range("K5").Value = "Panipat"

$ws.Range("A6").Value = "Ashwini"
$ws.Range("B6").Value = "Patel"
$ws.Range("D6").Value = "Female"
$ws.Range("E6").Value = 9876543210
$ws.Range("F6").Value = 34465
$ws.Range("G6").Value = "TAMIL, ENGLISH, MATHS, PHYSICS"
$ws.Range("H6").Value = "RADING"
$ws.Range("I6").Value = "PANIPAT"
$ws.Range("J6").Value = "Haryana"
$ws.Range("K6").Value = "Panipat"

# ---------------------------------------------------------------------------
# 4. EMAIL column values + hyperlinks for every data row. Set the display
#    text first so Hyperlinks.Add doesn't stamp a redundant display=""
#    attribute onto the <hyperlink> element.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "ranjith@gmail.com"
$ws.Range("C3").Value = "karthick@gmail.com"
$ws.Range("C4").Value = "saransant@gmail.com"
$ws.Range("C5").Value = "parthiban@gmail.com"
$ws.Range("C6").Value = "ashpat@gmail.com"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ranjith@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:karthick@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:saransant@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:parthiban@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:ashpat@gmail.com")

# ---------------------------------------------------------------------------
# 5. Cosmetics: selection.
# ---------------------------------------------------------------------------
$ws.Range("E10").Select() | Out-Null
